$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '62.461.86'
$ws.Cells.Item(2, 5).Value = '  -1.84%  '
$ws.Cells.Item(3, 4).Value = '2.438.52'
$ws.Cells.Item(3, 5).Value = '  -1.52%  '
$ws.Cells.Item(4, 5).Value = '  -0.23%  '
$ws.Cells.Item(5, 4).Value = "'565.29"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  -1.88%  '
$ws.Cells.Item(6, 4).Value = "'144.06"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  -3.29%  '
$ws.Cells.Item(7, 5).Value = '  +0.17%  '
$ws.Cells.Item(8, 4).Value = "'0.531"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -1.99%  '
$ws.Cells.Item(9, 4).Value = '2.434.57'
$ws.Cells.Item(9, 5).Value = '  -2.01%  '
$ws.Cells.Item(10, 5).Value = '  -5.61%  '
$ws.Cells.Item(11, 5).Value = '  +1.03%  '
$ws.Cells.Item(12, 4).Value = "'5.19"
$ws.Cells.Item(12, 4).Style = 'Normal'
$ws.Cells.Item(12, 5).Value = '  -3.04%  '
$ws.Cells.Item(13, 4).Value = "'0.353"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  -2.91%  '
$ws.Cells.Item(14, 4).Value = "'26.55"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  -6.01%  '
$ws.Cells.Item(16, 4).Value = '2.876.43'
$ws.Cells.Item(16, 5).Value = '  -1.69%  '
$ws.Cells.Item(17, 4).Value = '62.318.62'
$ws.Cells.Item(17, 5).Value = '  -1.90%  '
$ws.Cells.Item(18, 4).Value = '2.426.32'
$ws.Cells.Item(18, 5).Value = '  -2.40%  '
$ws.Cells.Item(19, 4).Value = "'11.12"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -4.03%  '
$ws.Cells.Item(20, 5).Value = '  -2.05%  '
$ws.Cells.Item(21, 4).Value = "'324.71"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  -1.30%  '
$ws.Cells.Item(22, 5).Value = '  -2.44%  '
$ws.Cells.Item(23, 4).Value = "'2.02"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  +5.68%  '
$ws.Cells.Item(24, 5).Value = '  +0.28%  '
$ws.Cells.Item(25, 4).Value = "'65.08"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  -3.79%  '
$ws.Cells.Item(26, 4).Value = "'625.08"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -0.95%  '
$ws.Cells.Item(27, 4).Value = "'9.00"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +2.38%  '
$ws.Cells.Item(28, 4).Value = '0.0₃0966'
$ws.Cells.Item(28, 5).Value = '  -9.03%  '
$ws.Cells.Item(29, 5).Value = '  -1.79%  '
$ws.Cells.Item(30, 5).Value = '  +0.34%  '
$ws.Cells.Item(31, 5).Value = '  -4.19%  '
$ws.Cells.Item(32, 4).Value = "'8.08"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  -4.66%  '
$ws.Cells.Item(33, 4).Value = "'1.87"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value = '  -2.25%  '
$ws.Cells.Item(34, 5).Value = '  -7.58%  '
$ws.Cells.Item(35, 4).Value = "'5.03"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  -3.74%  '
$ws.Cells.Item(36, 4).Value = "'1.00"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  +0.23%  '
$ws.Cells.Item(37, 4).Value = "'1.47"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -5.93%  '
$ws.Cells.Item(38, 5).Value = '  -2.90%  '
$ws.Cells.Item(39, 4).Value = "'18.75"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value = '  -1.45%  '
$ws.Cells.Item(40, 4).Value = "'5.26"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  -5.73%  '
$ws.Cells.Item(41, 4).Value = "'146.04"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  -0.48%  '
$ws.Cells.Item(42, 5).Value = '  -6.98%  '
$ws.Cells.Item(43, 2).Value = 'USDe'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(43, 4).Value = "'0.999"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  +0.00%  '
$ws.Cells.Item(44, 2).Value = 'OKB'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(44, 4).Value = "'42.18"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  +0.63%  '
$ws.Cells.Item(45, 2).Value = 'dogwifhat'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(45, 4).Value = "'2.53"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value = '  -5.10%  '
$ws.Cells.Item(46, 4).Value = "'144.99"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -4.15%  '
$ws.Cells.Item(47, 4).Value = "'3.71"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -1.96%  '
$ws.Cells.Item(48, 4).Value = "'20.31"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -4.50%  '
$ws.Cells.Item(49, 4).Value = "'0.0525"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -4.63%  '
$ws.Cells.Item(50, 4).Value = "'0.594"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -2.76%  '
$ws.Cells.Item(51, 5).Value = '  -4.60%  '
